# Auto-generated Excel COM-interop script
# Applies the 2024-06-13 daily violent-crime data update across all affected sheets.
# For every touched cell we just set the new absolute total (these sheets store
# running year-to-date totals per crime category / neighborhood, not formulas).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 3422
$ws.Range("K3").Value = 3390
$ws.Range("I4").Value = 1794
$ws.Range("K4").Value = 709
$ws.Range("K6").Value = 3997
$ws.Range("I7").Value = 26248
$ws.Range("K7").Value = 11740

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 80
$ws.Range("K7").Value = 156

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 229
$ws.Range("K3").Value = 235
$ws.Range("K7").Value = 783

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 91
$ws.Range("K7").Value = 252

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 174
$ws.Range("K6").Value = 135
$ws.Range("K7").Value = 471

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 141
$ws.Range("K7").Value = 405

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K6").Value = 112
$ws.Range("K7").Value = 284

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K6").Value = 90
$ws.Range("K7").Value = 334
$ws.Range("K8").Value = 783
$ws.Range("K10").Value = 65
$ws.Range("K11").Value = 244
$ws.Range("K14").Value = 61
$ws.Range("K17").Value = 21
$ws.Range("K19").Value = 365
$ws.Range("K20").Value = 270
$ws.Range("K22").Value = 37
$ws.Range("K23").Value = 115
$ws.Range("K25").Value = 49
$ws.Range("K27").Value = 120
$ws.Range("K29").Value = 619
$ws.Range("K31").Value = 126
$ws.Range("K33").Value = 471
$ws.Range("K36").Value = 139
$ws.Range("K37").Value = 405
$ws.Range("K44").Value = 109
$ws.Range("K48").Value = 144
$ws.Range("K49").Value = 69
$ws.Range("K52").Value = 318
$ws.Range("K53").Value = 156
$ws.Range("K55").Value = 125
$ws.Range("K57").Value = 36
$ws.Range("I63").Value = 209
$ws.Range("K63").Value = 41
$ws.Range("K65").Value = 284
$ws.Range("K67").Value = 455
$ws.Range("K71").Value = 35
$ws.Range("K73").Value = 105
$ws.Range("K76").Value = 179
$ws.Range("K78").Value = 150
$ws.Range("K79").Value = 303
$ws.Range("K83").Value = 252
$ws.Range("K89").Value = 157
$ws.Range("I101").Value = 26248
$ws.Range("K101").Value = 11740

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 126

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 141
$ws.Range("K7").Value = 455

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 69

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 173
$ws.Range("K3").Value = 208
$ws.Range("K6").Value = 188
$ws.Range("K7").Value = 619

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K4").Value = 20
$ws.Range("K7").Value = 144

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 124
$ws.Range("K3").Value = 97
$ws.Range("K6").Value = 116
$ws.Range("K7").Value = 365

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 179

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K2").Value = 36
$ws.Range("K7").Value = 90

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K3").Value = 33
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 30
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 38
$ws.Range("K3").Value = 40
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 303

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 77
$ws.Range("K7").Value = 270

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 54
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 139

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 123
$ws.Range("K3").Value = 100
$ws.Range("K7").Value = 334

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 74
$ws.Range("K7").Value = 244

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K4").Value = 4
$ws.Range("K6").Value = 63

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 120

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 192
$ws.Range("K6").Value = 127

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 84
$ws.Range("K7").Value = 318
